{"js": "// Update the \"Maximal (%)\" and \"Most recent (%)\" drop figures in the\n// vessel-category table. Every occurrence of each old value is replaced\n// with its corresponding new value (one pair repeats twice in the table).\nconst replacements = [\n  [\"-33.6\", \"-41.0\"],\n  [\"-14.8\", \"-16.0\"],\n  [\"-21.6\", \"-24.3\"],\n  [\"-14.2\", \"-15.3\"],\n  [\"-76.5\", \"-144.7\"],\n  [\"-54.0\", \"-77.8\"],\n  [\"-72.1\", \"-127.5\"],\n  [\"-12.8\", \"-13.7\"],\n  [\"-93.7\", \"-276.9\"],\n  [\"-28.8\", \"-34.0\"],\n  [\"-47.4\", \"-64.2\"],\n  [\"-62.3\", \"-97.5\"],\n  [\"-26.6\", \"-30.9\"],\n];\n\nfor (const [oldVal, newVal] of replacements) {\n  const results = context.document.body.search(oldVal, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newVal, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the \"Maximal (%)\" and \"Most recent (%)\" drop figures in the\n# vessel-category table (column 2 = Maximal %, column 4 = Most recent %,\n# data rows 3-9).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    3 = @{ 2 = \"-41.0\";  4 = \"-16.0\" }   # Cargo\n    4 = @{ 2 = \"-24.3\";  4 = \"-15.3\" }   # Tanker\n    5 = @{ 2 = \"-144.7\"; 4 = \"-77.8\" }   # Passenger\n    6 = @{ 2 = \"-127.5\"; 4 = \"-13.7\" }   # Fishing\n    7 = @{ 2 = \"-276.9\"; 4 = \"-34.0\" }   # Recreational\n    8 = @{ 2 = \"-64.2\";  4 = \"-16.0\" }   # Other\n    9 = @{ 2 = \"-97.5\";  4 = \"-30.9\" }   # All vessels\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $cols = $newValues[$rowIndex]\n    foreach ($colIndex in $cols.Keys) {\n        $cell = $t.Cell($rowIndex, $colIndex)\n        $r = $cell.Range\n        # Drop the trailing end-of-cell mark so only the visible text is replaced,\n        # preserving the run/paragraph formatting already on that text.\n        $r.End = $r.End - 1\n        $r.Text = $cols[$colIndex]\n    }\n}\n"}
